$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the placeholder "Promotional set" names (A/B/C) with real descriptive names
$ws.Range("J2").Value = "ASet House's special"
$ws.Range("J3").Value = "Set T-bone's special"
$ws.Range("J4").Value = "Set Fishy meal"

# Update the active cell selection
$ws.Range("K9").Select()
